$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 359, shifting the existing rows 359-448 down to 361-450
$ws.Rows.Item(359).Insert()
$ws.Rows.Item(359).Insert()

# New row 359 - Calidad "Primera"
$ws.Range("A359").Value = 8
$ws.Range("B359").Value = "Terminal La Palmera de La Serena"
$ws.Range("C359").Value = "Coquimbo"
$ws.Range("D359").Value = 44754
$ws.Range("E359").Value = 4
$ws.Range("F359").Value = 100112043
$ws.Range("G359").Value = "Pepino dulce"
$ws.Range("H359").Value = "Cultivar IV Región"
$ws.Range("I359").Value = "Primera"
$ws.Range("J359").Value = 400
$ws.Range("K359").Value = 13000
$ws.Range("L359").Value = 14000
$ws.Range("M359").Value = 13500
$ws.Range("N359").Value = '$/bandeja 18 kilos'
$ws.Range("O359").Value = "Provincia de Limarí"
$ws.Range("P359").Value = 750
$ws.Range("Q359").Value = 18
$ws.Range("R359").Value = "Hortaliza"

# New row 360 - Calidad "Segunda"
$ws.Range("A360").Value = 8
$ws.Range("B360").Value = "Terminal La Palmera de La Serena"
$ws.Range("C360").Value = "Coquimbo"
$ws.Range("D360").Value = 44754
$ws.Range("E360").Value = 4
$ws.Range("F360").Value = 100112043
$ws.Range("G360").Value = "Pepino dulce"
$ws.Range("H360").Value = "Cultivar IV Región"
$ws.Range("I360").Value = "Segunda"
$ws.Range("J360").Value = 280
$ws.Range("K360").Value = 10000
$ws.Range("L360").Value = 11000
$ws.Range("M360").Value = 10500
$ws.Range("N360").Value = "$/bandeja 18 kilos"
$ws.Range("O360").Value = "Provincia de Limarí"
$ws.Range("P360").Value = 583
$ws.Range("Q360").Value = 18
$ws.Range("R360").Value = "Hortaliza"
